# RPA datasets push 2024-05-30
# Insert a new IPO row ("신한글로벌액티브") above the "한중엔시에스" row,
# shifting the following rows down, and drop the former last row
# ("라메디텍") that falls off the bottom of the A1:F21 table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 18 ("한중엔시에스") and everything below it down by one row.
$ws.Rows.Item(18).Insert()

# Fill in the newly inserted row 18 with the new IPO's data.
$ws.Range("A18").Value = "신한글로벌액티브"
$ws.Range("B18").Value = "2024.06.03~06.05"
$ws.Range("C18").Value = "3,000~3,800"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = 70000
$ws.Range("F18").Value = "신한투자증권,한국투자증권"

# The table keeps its original extent (A1:F21), so the row that used to be
# last ("라메디텍", now shifted to row 22) is removed.
$ws.Rows.Item(22).Delete()
